$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 438, shifting all rows from
# 438 onward down by 3 (old 438 -> new 441, ... old 518 -> new 521).
$ws.Rows("438:440").Insert()

# New row 438: Tomate / Larga vida / Primera, fecha 2021-10-07 (44476)
$ws.Cells.Item(438,1).Value = 1
$ws.Cells.Item(438,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(438,3).Value = "Arica y Parinacota"
$ws.Cells.Item(438,4).Value = 44476
$ws.Cells.Item(438,5).Value = 15
$ws.Cells.Item(438,6).Value = 100112020
$ws.Cells.Item(438,7).Value = "Tomate"
$ws.Cells.Item(438,8).Value = "Larga vida"
$ws.Cells.Item(438,9).Value = "Primera"
$ws.Cells.Item(438,10).Value = 300
$ws.Cells.Item(438,11).Value = 5000
$ws.Cells.Item(438,12).Value = 5500
$ws.Cells.Item(438,13).Value = 5250
$ws.Cells.Item(438,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(438,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(438,16).Value = 525
$ws.Cells.Item(438,17).Value = 10
$ws.Cells.Item(438,18).Value = "Hortaliza"

# New row 439: Tomate / Larga vida / Segunda, fecha 2021-10-07 (44476)
$ws.Cells.Item(439,1).Value = 1
$ws.Cells.Item(439,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(439,3).Value = "Arica y Parinacota"
$ws.Cells.Item(439,4).Value = 44476
$ws.Cells.Item(439,5).Value = 15
$ws.Cells.Item(439,6).Value = 100112020
$ws.Cells.Item(439,7).Value = "Tomate"
$ws.Cells.Item(439,8).Value = "Larga vida"
$ws.Cells.Item(439,9).Value = "Segunda"
$ws.Cells.Item(439,10).Value = 350
$ws.Cells.Item(439,11).Value = 4000
$ws.Cells.Item(439,12).Value = 4500
$ws.Cells.Item(439,13).Value = 4250
$ws.Cells.Item(439,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(439,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(439,16).Value = 425
$ws.Cells.Item(439,17).Value = 10
$ws.Cells.Item(439,18).Value = "Hortaliza"

# New row 440: Tomate / Larga vida / Tercera, fecha 2021-10-07 (44476)
$ws.Cells.Item(440,1).Value = 1
$ws.Cells.Item(440,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(440,3).Value = "Arica y Parinacota"
$ws.Cells.Item(440,4).Value = 44476
$ws.Cells.Item(440,5).Value = 15
$ws.Cells.Item(440,6).Value = 100112020
$ws.Cells.Item(440,7).Value = "Tomate"
$ws.Cells.Item(440,8).Value = "Larga vida"
$ws.Cells.Item(440,9).Value = "Tercera"
$ws.Cells.Item(440,10).Value = 400
$ws.Cells.Item(440,11).Value = 3500
$ws.Cells.Item(440,12).Value = 4000
$ws.Cells.Item(440,13).Value = 3750
$ws.Cells.Item(440,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(440,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(440,16).Value = 375
$ws.Cells.Item(440,17).Value = 10
$ws.Cells.Item(440,18).Value = "Hortaliza"
